$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.208.23"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.01"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6998"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.63"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07888"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3016"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.94"
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08136"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.858.01"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.190"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7065"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.39"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.196.90"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.808"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007832"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.62"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.089.86"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.496"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.54"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.882"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1417"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.917"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.471"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.302"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.010"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.167"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7063"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9962"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.705"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.151.84"
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9210"
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.946"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4238"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.0000"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.97"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5292"
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.736"
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.159"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.965"
$ws.Range("E51").Value = "  -0.30%  "
